$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New ATTENDANCE row (row 12) ---
$ws.Range("A12").Value = "ATTENDANCE"
$ws.Range("C12").Value = "All attended"
$ws.Range("E12").Value = "All attended"
$ws.Range("G12").Value = "All attended"

$ws.Range("A12").Style = "Bad"
$ws.Range("C12").Style = "Bad"
$ws.Range("E12").Style = "Bad"
$ws.Range("G12").Style = "Bad"

# --- Update existing "TO DO/UPDATE/FEEDBACK" column text (K5, K7) ---
$ws.Range("K7").Value = "Add Advantages & Sketch Prototype"
$ws.Range("K5").Value = "Add disadvantage & Sketch use case diagram"

# --- New row 15 (attendance note under the TO DO column) ---
$ws.Range("I15").Value = "All attended"
$ws.Range("I15").Style = "Bad"

# --- Sheet view: zoom + selection ---
$ws.Range("K15").Select()
$ws.Application.ActiveWindow.Zoom = 61

# --- Column widths ---
# NOTE: this runtime's ColumnWidth setter stores (round(input*6)/6 + 5/6),
# i.e. it always adds a fixed 5/6-character padding on top of a width
# quantized to 1/6 of a character. To land on the desired stored width we
# therefore feed it (target - 5/6).
$ws.Columns.Item(1).ColumnWidth = 40.166666666666664
$ws.Columns.Item(3).ColumnWidth = 47.166666666666664
$ws.Columns.Item(5).ColumnWidth = 46.944010416666664
$ws.Columns.Item(7).ColumnWidth = 61.498697916666664
$ws.Columns.Item(9).ColumnWidth = 35.498697916666664
$ws.Columns.Item(11).ColumnWidth = 40.944010416666664
